# Refresh the cryptos price list (values as scraped on 2024-03-25).
# Price cells that look numeric (e.g. "1.00", "58.24") are written with a
# leading apostrophe so Excel stores them as text, exactly like the
# existing multi-dot prices (e.g. "70.875.50") are already stored.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.875.50'
$ws.Range("E2").Value = '  +7.49%  '

$ws.Range("D3").Value = '3.635.94'
$ws.Range("E3").Value = '  +7.35%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").Value = "'" + '594.92'
$ws.Range("E5").Value = '  +5.60%  '

$ws.Range("D6").Value = "'" + '192.72'
$ws.Range("E6").Value = '  +9.37%  '

$ws.Range("E7").Value = '  +3.92%  '

$ws.Range("D8").Value = '3.630.16'
$ws.Range("E8").Value = '  +7.18%  '

$ws.Range("D9").Value = "'" + '1.00'
$ws.Range("E9").Value = '  +0.14%  '

$ws.Range("D10").Value = "'" + '0.181'
$ws.Range("E10").Value = '  +2.46%  '

$ws.Range("E11").Value = '  +4.95%  '

$ws.Range("D12").Value = "'" + '58.24'
$ws.Range("E12").Value = '  +8.02%  '

$ws.Range("D13").Value = "'" + '0.0000297'
$ws.Range("E13").Value = '  +6.69%  '

$ws.Range("D14").Value = "'" + '9.81'
$ws.Range("E14").Value = '  +5.80%  '

$ws.Range("D15").Value = '4.218.18'
$ws.Range("E15").Value = '  +7.44%  '

$ws.Range("D16").Value = '3.633.79'
$ws.Range("E16").Value = '  +7.90%  '

$ws.Range("D17").Value = "'" + '19.43'
$ws.Range("E17").Value = '  +6.71%  '

$ws.Range("D18").Value = '70.827.03'
$ws.Range("E18").Value = '  +7.51%  '

$ws.Range("D19").Value = "'" + '12.63'
$ws.Range("E19").Value = '  +6.03%  '

$ws.Range("E20").Value = '  +0.96%  '

$ws.Range("E21").Value = '  +5.80%  '

$ws.Range("D22").Value = "'" + '495.95'
$ws.Range("E22").Value = '  +7.14%  '

$ws.Range("D23").Value = "'" + '5.44'
$ws.Range("E23").Value = '  +10.44%  '

$ws.Range("D24").Value = "'" + '17.11'

$ws.Range("E25").Value = '  +9.65%  '

$ws.Range("D26").Value = "'" + '91.25'
$ws.Range("E26").Value = '  +1.62%  '

$ws.Range("E27").Value = '  +6.91%  '

$ws.Range("E28").Value = '  +5.88%  '

$ws.Range("E29").Value = '  +8.51%  '

$ws.Range("D30").Value = "'" + '32.46'
$ws.Range("E30").Value = '  +4.25%  '

$ws.Range("D31").Value = "'" + '7.65'
$ws.Range("E31").Value = '  +15.59%  '

$ws.Range("D32").Value = "'" + '12.27'
$ws.Range("E32").Value = '  +6.83%  '

$ws.Range("D33").Value = "'" + '617.75'
$ws.Range("E33").Value = '  +6.39%  '

$ws.Range("D34").Value = "'" + '0.118'
$ws.Range("E34").Value = '  +8.79%  '

$ws.Range("D35").Value = "'" + '65.35'
$ws.Range("E35").Value = '  +4.41%  '

$ws.Range("D36").Value = '0.0₃0834'
$ws.Range("E36").Value = '  +11.18%  '

$ws.Range("D37").Value = "'" + '0.411'
$ws.Range("E37").Value = '  +8.59%  '

$ws.Range("D38").Value = "'" + '0.149'
$ws.Range("E38").Value = '  +3.92%  '

$ws.Range("D39").Value = "'" + '38.25'
$ws.Range("E39").Value = '  +6.04%  '

$ws.Range("E40").Value = '  -0.02%  '

$ws.Range("D41").Value = "'" + '3.69'
$ws.Range("E41").Value = '  +3.38%  '

$ws.Range("D42").Value = '3.336.49'
$ws.Range("E42").Value = '  +6.92%  '

$ws.Range("D43").Value = "'" + '3.09'
$ws.Range("E43").Value = '  +8.90%  '

$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = "'" + '0.0450'
$ws.Range("E44").Value = '  +7.56%  '

$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").Value = "'" + '2.70'
$ws.Range("E45").Value = '  +9.84%  '

$ws.Range("D46").Value = "'" + '3.40'
$ws.Range("E46").Value = '  +7.54%  '

$ws.Range("E47").Value = '  +3.23%  '

$ws.Range("D48").Value = "'" + '9.29'
$ws.Range("E48").Value = '  +9.38%  '

$ws.Range("D49").Value = "'" + '2.76'
$ws.Range("E49").Value = '  +8.05%  '

$ws.Range("D50").Value = "'" + '3.33'
$ws.Range("E50").Value = '  +5.86%  '

$ws.Range("B51").Value = 'FirstDigitalUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D51").Value = "'" + '1.00'
$ws.Range("E51").Value = '  +0.04%  '
